$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New vocabulary added to "Lektion 2: Haus und Kleidung" (columns F/G)
$ws.Range("F31").Value = "Teppich"
$ws.Range("G31").Value = "carpet"
$ws.Range("F32").Value = "Badewanne"
$ws.Range("G32").Value = "bathtub"
$ws.Range("F33").Value = "Regal"
$ws.Range("G33").Value = "shelf"
$ws.Range("F34").Value = "Pullover"
$ws.Range("G34").Value = "pullover"
$ws.Range("F35").Value = "Mantel"
$ws.Range("G35").Value = "coat"

# New vocabulary added to "Lektion 3: Tiere" (columns J/K)
$ws.Range("J30").Value = "Tiger"
$ws.Range("K30").Value = "tiger"
$ws.Range("J31").Value = "Delfin"
$ws.Range("K31").Value = "dolphin"
$ws.Range("J32").Value = "Pinguin"
$ws.Range("K32").Value = "penguin"
$ws.Range("J33").Value = "Maus"
$ws.Range("K33").Value = "mouse"
$ws.Range("J34").Value = "Kuh"
$ws.Range("K34").Value = "cow"
$ws.Range("J35").Value = "Ratte"
$ws.Range("K35").Value = "rat"
$ws.Range("J36").Value = "Schaf"
$ws.Range("K36").Value = "sheep"
$ws.Range("J37").Value = "Ente"
$ws.Range("K37").Value = "duck"
$ws.Range("J38").Value = "Gans"
$ws.Range("K38").Value = "goose"

# Update view state to match the author's final selection position
$ws.Range("E36").Select()
